$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.585.48'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '2.638.58'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.10%  '
$ws.Range('E7').Value = '  +2.60%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -2.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.84'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.390'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.85'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000188'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.56%  '
$ws.Range('D15').Value = '3.120.69'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').Value = '64.317.46'
$ws.Range('E16').Value = '  -1.64%  '
$ws.Range('D17').Value = '2.648.19'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.26'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.71'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.53'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '348.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.82%  '
$ws.Range('E24').Value = '  +9.15%  '
$ws.Range('E25').Value = '  +0.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.47'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '595.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +10.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.04'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.50%  '
$ws.Range('E30').Value = '  -0.88%  '
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('E32').Value = '  -0.56%  '
$ws.Range('E33').Value = '  +0.63%  '
$ws.Range('E34').Value = '  +4.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.35'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.40%  '
$ws.Range('E36').Value = '  -0.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.09'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('E39').Value = '  +2.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '152.82'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.09%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '159.19'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('E43').Value = '  +5.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.03'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '23.55'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.57%  '
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('E47').Value = '  +0.54%  '
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('E49').Value = '  +2.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.27'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.18%  '
$ws.Range('D51').Value = '0.0₆0236'
$ws.Range('E51').Value = '  -5.86%  '
